# Add 8 new Spanish/English vocabulary rows (133-140) to the "theenglish"
# table/sheet, clear the leftover "applyNumberFormat" style from the
# previously-last rows (122-132) so it becomes unused again, extend the
# DatosExternos_1 defined name to match the new table size, and update the
# sheet view / selection to reflect where the user ended up scrolled to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("theenglish")

# Rows 122-132 previously carried an extra (now pointless) number-format
# style; drop it so the cells go back to the default style.
$ws.Range("A122:C132").ClearFormats()

# New rows to append to the "theenglish" table (esp / eng / structure).
$newRows = @(
    @("No lleves cosas pesadas", "don't carry heavy things", "all"),
    @("Pide a laguien que te ayude", "Ask someone to help you", "all"),
    @("Ya no trbajan aquí", "They don't work here anymore", "all"),
    @("Recuerda lo que pasó la última vez", "Remember what happened last time.", "all"),
    @("Creia que estaban hablando en Alemán", "I thought they were speaking German", "all"),
    @("He oido suficiente", "I've heard enought", "all"),
    @("El le mas rádido que la mayoría de la gente", "He reads faster than most people", "all"),
    @("El programa durará hasta la próxima semana", "The program will last until next week", "all")
)

$lo = $ws.ListObjects.Item(1)

foreach ($rowData in $newRows) {
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range.Row
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

# Keep the external-data defined name range in sync with the table size.
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
$nm = $wb.Names.Item("DatosExternos_1")
$nm.RefersTo = "=theenglish!`$A`$1:`$B`$" + $lastRow

# Reflect the scrolled/selected position left behind by the edit.
$ws.Range("B132").Select()
$excel.ActiveWindow.ScrollRow = 109
